$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the model strings in column A (values unchanged order, renamed per new naming convention)
$ws.Range("A2").Value = "year + site"
$ws.Range("A3").Value = "coralcover + year + site"
$ws.Range("A4").Value = "year + site + year*site"
$ws.Range("A5").Value = "site"
$ws.Range("A6").Value = "coralcover + site"
$ws.Range("A7").Value = "coralcover + site + site*coralcover"
$ws.Range("A8").Value = "coralcover + year"
$ws.Range("A9").Value = "coralcover + year + year*coralcover"
$ws.Range("A10").Value = "coralcover"
$ws.Range("A11").Value = "year"

# Widen column A to fit the new, longer text (COM ColumnWidth maps to a
# stored XML width that is ~0.8333 wider, so back that padding out to land
# on an exact width of 34)
$ws.Columns.Item(1).ColumnWidth = 33.16666666666666

# Reset the selection stored in the sheet view back to the default top-left cell
$ws.Range("A1").Select()
